$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 'Insperity'
$ws.Range("B7").Value = 'This system is primarily used for HR management and employee administration. It supports onboarding, I-9 verification, benefits election, payroll processing, time reporting, vacation accruals, and time tracking.'
$ws.Range("C7").Value = 'The client utilizes this system to manage onboarding processes, employee benefit elections, time tracking for vacation accruals, payroll processing, and manual compensation updates within employee profiles.'
$ws.Range("D7").Value = 'The HR team, led by Laurie Stewart (SVP of HR), manages system administration and user access.'
$ws.Range("E7").Value = 'Access provisioning is initiated when HR inputs details such as salary, start date, and manager information for new hires, triggering a notification for onboarding paperwork and granting employee access to Insperity. Elevated access requires an admin to submit a permission-based access request form to the Insperity team for provisioning, as Jade users cannot directly provision access. Role changes are not explicitly detailed but likely follow the elevated access request process.'
$ws.Range("F7").Value = 'Upon termination, HR updates the employee''s termination date and reason in the system. For elevated access users, an email is sent to Insperity to request access removal, and confirmation is received via email documenting the process.'
$ws.Range("G7").Value = 'Access is configured using a permission-based model, where specific permissions are selected through an access request form and provisioned by the Insperity team.'
$ws.Range("H7").Value = 'No; Access is permission-based, and system administrators cannot create, modify, or delete roles. Changes to permissions require submission of an access request form processed by the vendor.'
$ws.Range("I7").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("J7").Value = 'Yes; Users with privileged access include Megan Hodgson, John, and the HR team, who can perform unrestricted functions such as setting up new employees, though payroll processing is controlled externally by Insperity.'
$ws.Range("K7").Value = 'No; there are no interactive system accounts within the system.'
$ws.Range("L7").Value = 'The credentials for all accounts are fully managed and stored by Insperity''s system; there is no local storage or management of these credentials by the Jade team.'
$ws.Range("M7").Value = 'N/A - No Interactive System Accounts'
$ws.Range("N7").Value = 'Yes; Management will perform periodic user access reviews. No additional detail was provided regarding the frequency or process during the walkthrough meeting.'
$ws.Range("O7").Value = 'The system maintains logs of pay change activities. Logs include change history details but do not explicitly confirm comprehensive admin activity logging or audit trail functionality.'
$ws.Range("P7").Value = 'No; periodic activity reviews of user activity are not currently performed, but there are related controls such as payroll reviews and plans to implement user access reviews in the future.'
$ws.Range("Q7").Value = 'Users authenticate through direct login with enforced two-factor authentication for all accounts.'
$ws.Range("R7").Value = 'No; The client does not perform periodic reviews of the systems authentication configurations.'
$ws.Range("S7").Value = 'Management does not have the ability to make any changes to this system. All system configurations, workflows, and code are managed entirely by the vendor.'
$ws.Range("T7").Value = 'Access to make changes is restricted to Megan Hodgson, John, and the HR team, who have administrative access to the application. Megan and John have admin access due to the current size of the team and operational needs, while the HR team has admin access for managing system configurations. All change access is granted through role-based permissions in the system.'
$ws.Range("U7").Value = 'Management does not have any separate environments for this system.'
$ws.Range("V7").Value = 'The change management process varies by change type. Access changes require an access request form specifying permissions, internal approval, and provisioning by Insperity, with Jade employees restricted from direct modifications. Payroll changes are requested via email, reviewed internally by designated personnel, and processed by Insperity with draft reviews for approval. Termination of user access involves email requests to Insperity, documented approvals, and confirmation of access disablement. Periodic payroll reviews serve as indirect validation of payroll changes, and a pay change history report tracks payroll-related modifications. There is no formalized process for testing and development, but plans to implement a new payroll system suggest future formalization of change management controls.'
$ws.Range("W7").Value = 'No, the vendor manages all updates, patches, and bug fixes directly in the production environment. The client does not have a sandbox or QA environment, and no changes are made internally by the client''s IT team.'
$ws.Range("X7").Value = 'N/A - This information was not discussed in the walkthrough meeting transcript.'
$ws.Range("Y7").Value = 'No periodic review of changes is performed.'
$ws.Range("Z7").Value = 'No automated jobs or interfaces are currently implemented for this system.'
$ws.Range("AA7").Value = 'The tools used to run, schedule, and monitor the automated jobs were not identified, and no specific capabilities were discussed.'
$ws.Range("AB7").Value = 'The failure resolution process involves contacting Insperity''s support team, waiting for them to resolve the issue, and relying on vendor-managed operations for backups and system-related issues.'
$ws.Range("AC7").Value = 'Data is stored in vendor-managed systems handled by Insperity.'
$ws.Range("AD7").Value = 'Backups are handled entirely by the vendor as part of their SaaS service.'
$ws.Range("AE7").Value = 'N/A - All backup types and strategies are managed by the vendor as part of their SaaS service.'
$ws.Range("AF7").Value = 'N/A - the vendor is responsible for monitoring and resolving any backup failures.'
$ws.Range("AG7").Value = 'No, management does not perform regular SOC report reviews.'
